$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# 1) Workbook calc setting: switch reference style to R1C1
$excel.ReferenceStyle = [Microsoft.Office.Interop.Excel.XlReferenceStyle]::xlR1C1

# 2) Fix two shared-string typos (Makefile -> makefile) used in D2/D3
$ws.Range("D2").Value = "Revision of makefile"
$ws.Range("D3").Value = "Documentation makefile changes. Concept for new sync objects"

# 3) Row 28 gets an additional-effort value in column C
$ws.Range("C28").Value = 0.25

# 4) Append new row 30 continuing the "Implementation tc14" task
#    Copy the date cell's format from the row above so the new date cell
#    keeps the same style (and does not spawn a duplicate number format).
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("A30").Value = (Get-Date -Year 2013 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B30").Value = 2
$ws.Range("D30").Value = $ws.Range("D29").Value2

# 5) Update selection to reflect the new last cell
$ws.Range("D30").Select()
